$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: add Status value in F5
$ws.Range("F5").Value = "PASS"

# Row 6: new login entry
$ws.Range("C6").Value = "epljxtvbgwbcxyl@gmail.com"
$ws.Range("D6").Value = "iltovYZDUX5"
$ws.Range("E6").Value = "pass"
$ws.Range("F6").Value = "PASS"

# Row 7: new login entry (no Status value)
$ws.Range("C7").Value = "ljrhsjgymlygiyx@gmail.com"
$ws.Range("D7").Value = "rielwWEGCW5"
$ws.Range("E7").Value = "pass"
